$d = $word.ActiveDocument

$replacements = @(
    @("2024-05-29 Wednesday", "2024-05-30 Thursday"),
    @("920÷6=", "143÷9="),
    @("357÷4=", "823÷3="),
    @("485÷3=", "679÷4="),
    @("422÷2=", "683÷6="),
    @("882÷6=", "900÷2="),
    @("758÷9=", "509÷2="),
    @("128÷9=", "570÷3="),
    @("766÷3=", "657÷5="),
    @("232÷2=", "249÷3="),
    @("403÷8=", "655÷8="),
    @("886÷8=", "306÷2="),
    @("333÷5=", "391÷6="),
    @("363÷5=", "407÷8="),
    @("678÷3=", "764÷8="),
    @("329÷6=", "338÷6="),
    @("698÷5=", "960÷4="),
    @("115÷4=", "287÷6="),
    @("841÷6=", "225÷2="),
    @("648÷9=", "278÷2="),
    @("402÷2=", "620÷9="),
    @("190÷9=", "503÷6="),
    @("536÷5=", "541÷6="),
    @("866÷5=", "938÷2="),
    @("200÷6=", "272÷4="),
    @("684÷9=", "855÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
